$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column M (w_audit_usefulness) first so column indices of earlier columns
# (G) remain valid, then delete column G (audit_usefulness).
$ws.Range("M1").EntireColumn.Delete()
$ws.Range("G1").EntireColumn.Delete()

# Update row 2 values per the diff
$ws.Range("B2").Value = 3
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 0.75
$ws.Range("H2").Value = 0.75
$ws.Range("I2").Value = 0.4
$ws.Range("J2").Value = 0.3
$ws.Range("K2").Value = 0.6
$ws.Range("L2").Value = 2.8
$ws.Range("M2").Value = "The report provides a reasonable extraction of evidence with correct citations, but lacks full sentences in some excerpts, affecting evidence extraction quality. Coverage of debiasing methods is broad, yet lacks depth in validation details, particularly for pre-processing and post-processing methods. The structure is clear and well-organized, aiding readability. However, relevance is compromised by some unsupported assumptions, such as the validation status of certain methods. Missing disclosures are not thoroughly identified, especially regarding the absence of validation metrics. The audit usefulness is moderate, as the report is traceable but lacks detailed validation evidence, limiting its utility for comprehensive audits."
